# "Added Enquiry to remaining Params"
# Adds a "Coverage Name" / "Allowed Coverages" mini reference table in
# columns Z:AA of Sheet1, alongside the existing PARAM/COLUMN grid.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- column widths for the two new columns -------------------------------
$ws.Columns.Item(26).ColumnWidth = 46.166666666666664   # Z  -> stored width 47
$ws.Columns.Item(27).ColumnWidth = 41.0                  # AA -> stored width ~41.83 (closest)

# --- header pair: Z2 / AA2 (bordered, wrapped, small font) ---------------
$ws.Range("Z2").Value = "Coverage Name"
$ws.Range("AA2").Value = "Allowed Coverages"

# --- label column Z3:Z7 (bordered, wrapped, small font) -------------------
$ws.Range("Z3").Value = "Is the Coverage mandatory"
$ws.Range("Z4").Value = "Is this coverage Basisc or a Rider"
$ws.Range("Z5").Value = "Can the Term exceed the Term of Basic"
$ws.Range("Z6").Value = "Can Premium Paying Term exceed the Premium Paying Term of Basic Coverage"
$ws.Range("Z7").Value = "Can the Term of Rider exceed the Term of Basic Coverage"

# Apply the bordered "box" style (font size 7.5, vertical-center + wrap,
# thin box border on all sides) to the header row pair and the label column.
$bordered = $ws.Range("Z2,AA2,Z3,Z4,Z5,Z6,Z7")
$bordered.Font.Size = 7.5
$bordered.VerticalAlignment = -4108   # xlVAlignCenter
$bordered.WrapText = $true
$bordered.Borders.Weight = 2          # xlThin
$bordered.Borders.LineStyle = 1       # xlContinuous

# --- the matching (currently empty) AA3:AA7 cells --------------------------
# Same font/alignment, but no border - just reserved/blank cells waiting for
# values per-row.
$plain = $ws.Range("AA3,AA4,AA5,AA6,AA7")
$plain.Font.Size = 7.5
$plain.VerticalAlignment = -4108      # xlVAlignCenter
$plain.WrapText = $true

# Row 6 needs extra height because of the long wrapped label text.
$ws.Rows.Item(6).RowHeight = 20.399999999999999

# --- view state: scroll right and select the new table ---------------------
$ws.Range("Z2:AA7").Select()
